$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data-wide-value")

# Revert "User data 3.0": remove the "budget-type" column (column B) from the
# wide-value sheet, shifting columns C:G left into B:F.
$ws.Columns("B").Delete()
